$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing rows with revised International (B) / National (C) figures ---
# Row 2
$ws.Range("B2").Value = 235428
# Row 4
$ws.Range("B4").Value = 756970
# Row 6
$ws.Range("B6").Value = 1005152
# Row 7
$ws.Range("B7").Value = 867801
# Row 8
$ws.Range("B8").Value = 399342
$ws.Range("C8").Value = 2158520
# Row 9
$ws.Range("C9").Value = 2592733
# Row 10
$ws.Range("B10").Value = 515952
$ws.Range("C10").Value = 7194502
# Row 11
$ws.Range("B11").Value = 511365
$ws.Range("C11").Value = 3337189
# Row 12
$ws.Range("B12").Value = 1311766
# Row 13
$ws.Range("B13").Value = 1028704
$ws.Range("C13").Value = 1393812
# Row 14
$ws.Range("B14").Value = 797443
$ws.Range("C14").Value = 3417580

# --- Extend the table (Tabla1) with 6 new rows of data (rows 15-20) ---
$tbl = $ws.ListObjects.Item("Tabla1")
$tbl.Resize($ws.Range("A1:D20")) | Out-Null

$newRows = @(
    @(15, 43538, 780152, 1034335),
    @(16, 43539, 368021, 2250080),
    @(17, 43540, 633122, 4918930),
    @(18, 43541, 982658, 6175263),
    @(19, 43542, 2038317, 6219920),
    @(20, 43543, 1118360, 4845820)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $dateVal = $row[1]
    $intl = $row[2]
    $natl = $row[3]

    # Copy the date-formatted style from the row above so the new date cell
    # matches the existing column A formatting (numFmtId 14 / style index 1).
    $ws.Range("A$($r - 1)").Copy($ws.Range("A$r")) | Out-Null
    $ws.Range("A$r").Value = $dateVal

    $ws.Range("B$r").Value = $intl
    $ws.Range("C$r").Value = $natl
    $ws.Range("D$r").Formula = '=B' + $r + '+C' + $r + '/Hoja2!$A$2'
}

# --- UI state: active cell / selection on D8 ---
$ws.Range("D8").Select() | Out-Null
